{"js": "// Fixed duplicated values when computing institutional factors.\n// Update the p-values table (Fig 3 trend lines) in place. Several cells\n// share the same old value (e.g. two \"0.19\" cells, two \"<0.01\" cells), so\n// the lookup is scoped per-cell (by row label + column) and the text is\n// replaced via a cell-scoped search instead of a whole-document one, to\n// avoid touching the wrong occurrence and to keep each cell's existing\n// (empty + value) run pair intact.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// row label -> { columnIndex (1 = Cod, 2 = Hake) -> [oldValue, newValue] }\nconst updates = {\n  \"GDP 2016\": { 1: [\"0.56\", \"0.41\"], 2: [\"0.19\", \"0.17\"] },\n  \"OHI economic 2016\": { 1: [\"0.99\", \"0.79\"], 2: [\"0.89\", \"0.97\"] },\n  \"Technical Development (number per country)\": { 1: [\"0.24\", \"0.28\"], 2: [\"0.05\", \"0.06\"] },\n  \"Compilance (scores)\": { 1: [\"0.37\", \"0.43\"], 2: [\"0.25\", \"0.30\"] },\n  \"Readiness\": { 1: [\"0.15\", \"0.10\"], 2: [\"0.12\", \"0.11\"] },\n  \"Vulnerability\": { 1: [\"0.02\", \"0.01\"] },\n};\n\nconst values = table.values;\nconst pending = [];\nfor (let r = 0; r < values.length; r++) {\n  const rowLabel = values[r][0];\n  const rowUpdates = updates[rowLabel];\n  if (!rowUpdates) continue;\n  for (const colIndex of Object.keys(rowUpdates)) {\n    const c = Number(colIndex);\n    const [oldVal, newVal] = rowUpdates[colIndex];\n    if (values[r][c] !== oldVal) continue;\n    // Scope the search to this single cell so duplicate values elsewhere\n    // in the table (e.g. the other \"0.19\") are never touched.\n    const results = table.getCell(r, c).body.search(oldVal, { matchWholeWord: true });\n    results.load(\"items\");\n    pending.push({ results, newVal });\n  }\n}\nawait context.sync();\n\nfor (const { results, newVal } of pending) {\n  if (results.items.length !== 1) {\n    throw new Error(`expected exactly one match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(newVal, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fixed duplicated values when computing institutional factors.\n# Update the p-values table (Fig 3 trend lines) in place, cell by cell,\n# so that rows sharing the same old value (e.g. the two \"0.19\" / \"<0.01\"\n# cells) are only updated at the intended row/column.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# column 2 = Cod, column 3 = Hake (1-indexed; column 1 is the row label)\n$updates = @{\n    \"GDP 2016\"                                     = @{ 2 = @(\"0.56\", \"0.41\"); 3 = @(\"0.19\", \"0.17\") }\n    \"OHI economic 2016\"                             = @{ 2 = @(\"0.99\", \"0.79\"); 3 = @(\"0.89\", \"0.97\") }\n    \"Technical Development (number per country)\"    = @{ 2 = @(\"0.24\", \"0.28\"); 3 = @(\"0.05\", \"0.06\") }\n    \"Compilance (scores)\"                           = @{ 2 = @(\"0.37\", \"0.43\"); 3 = @(\"0.25\", \"0.30\") }\n    \"Readiness\"                                     = @{ 2 = @(\"0.15\", \"0.10\"); 3 = @(\"0.12\", \"0.11\") }\n    \"Vulnerability\"                                 = @{ 2 = @(\"0.02\", \"0.01\") }\n}\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowLabel = $t.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n    if (-not $updates.ContainsKey($rowLabel)) { continue }\n    $rowUpdates = $updates[$rowLabel]\n    foreach ($c in $rowUpdates.Keys) {\n        $pair = $rowUpdates[$c]\n        $oldVal = $pair[0]\n        $newVal = $pair[1]\n        $cellRange = $t.Cell($r, $c).Range\n        $curText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($curText -eq $oldVal) {\n            $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)\n            $textRange.Text = $newVal\n        }\n    }\n}\n"}
